$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '44.576.57'
$ws.Range("E2").Value = '  +0.65%  '

# Row 3
$ws.Range("D3").Value = '2.253.61'
$ws.Range("E3").Value = '  -0.41%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +0.27%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.48'
$ws.Range("E5").Value = '  +0.47%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.34'
$ws.Range("E6").Value = '  -3.01%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.574'
$ws.Range("E7").Value = '  -0.33%  '

# Row 8
$ws.Range("E8").Value = '  +0.35%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.529'
$ws.Range("E9").Value = '  -2.54%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.40'
$ws.Range("E10").Value = '  -2.27%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0815'
$ws.Range("E11").Value = '  -1.18%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.24'
$ws.Range("E12").Value = '  -2.47%  '

# Row 13
$ws.Range("E13").Value = '  +0.05%  '

# Row 14
$ws.Range("D14").Value = '2.366.12'
$ws.Range("E14").Value = '  +5.08%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.838'
$ws.Range("E15").Value = '  -1.18%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.66'
$ws.Range("E16").Value = '  -2.58%  '

# Row 17
$ws.Range("D17").Value = '44.362.50'
$ws.Range("E17").Value = '  +0.36%  '

# Row 18
$ws.Range("D18").Value = '0.0₃0967'
$ws.Range("E18").Value = '  -1.34%  '

# Row 19
$ws.Range("B19").Value = 'InternetComputer(DFINITY)'
$ws.Range("C19").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.41'
$ws.Range("E19").Value = '  -5.16%  '

# Row 20
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.37'
$ws.Range("E20").Value = '  -0.40%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '65.83'
$ws.Range("E21").Value = '  +0.18%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '238.66'
$ws.Range("E22").Value = '  -1.32%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.96'
$ws.Range("E23").Value = '  +0.40%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.02'
$ws.Range("E24").Value = '  +0.28%  '

# Row 25
$ws.Range("E25").Value = '  -0.19%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '40.39'
$ws.Range("E26").Value = '  +6.72%  '

# Row 27
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.92'
$ws.Range("E27").Value = '  -3.15%  '

# Row 28
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.21'
$ws.Range("E28").Value = '  +3.78%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.96'
$ws.Range("E29").Value = '  -2.83%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.16'
$ws.Range("E30").Value = '  -0.75%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '152.76'
$ws.Range("E31").Value = '  -2.81%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0804'
$ws.Range("E32").Value = '  -4.44%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.61'
$ws.Range("E33").Value = '  -2.32%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.15'
$ws.Range("E34").Value = '  -8.97%  '

# Row 35
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.111'
$ws.Range("E35").Value = '  +2.83%  '

# Row 36
$ws.Range("B36").Value = 'Stellar'
$ws.Range("C36").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.120'
$ws.Range("E36").Value = '  +0.61%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.79'
$ws.Range("E37").Value = '  -7.13%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.52'
$ws.Range("E38").Value = '  +0.45%  '

# Row 39
$ws.Range("B39").Value = 'Celestia'
$ws.Range("C39").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.78'
$ws.Range("E39").Value = '  -9.06%  '

# Row 40
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.87'
$ws.Range("E40").Value = '  -1.70%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0305'
$ws.Range("E41").Value = '  -1.45%  '

# Row 42
$ws.Range("E42").Value = '  +0.37%  '

# Row 43
$ws.Range("D43").Value = '1.734.01'
$ws.Range("E43").Value = '  -0.20%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '82.51'
$ws.Range("E44").Value = '  -7.83%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.195'
$ws.Range("E45").Value = '  -0.40%  '

# Row 46
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.64'
$ws.Range("E46").Value = '  +2.40%  '

# Row 47
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '100.23'
$ws.Range("E47").Value = '  -2.45%  '

# Row 48
$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.97'
$ws.Range("E48").Value = '  -4.25%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '55.72'
$ws.Range("E49").Value = '  -1.34%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.18'
$ws.Range("E50").Value = '  -1.92%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '69.78'
$ws.Range("E51").Value = '  -2.66%  '

